$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute(
    "Artificial Intelligence Transforms Medicine", $true, $false, $false, $false, $false,
    $true, 1, $false, "Unraveling the Enigmatic World of Chemistry", 2) | Out-Null

# --- Author line: three runs ("Dr" + "." + " Emily Kingston") collapse into one run ---
$r = $d.Content
$r.Find.Execute("Dr. Emily Kingston", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r.Text = "Brandy Harper"

# --- Email line: three runs ("ek@healthsciences" + "." + "org") collapse into one run ---
# Set .Text directly (instead of Find's Replace) so the straight apostrophe in
# "Author's" is not mangled by smart-quote autocorrect.
$r = $d.Content
$r.Find.Execute("ek@healthsciences.org", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$r.Text = "Author's Email address"

# --- Body paragraph (first block) ---
$d.Content.Find.Execute(
    "The convergence of artificial intelligence (AI) and medicine is paving the way for unprecedented advancements in healthcare",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Step into the captivating realm of chemistry, a science that holds the key to understanding the world around us", 2) | Out-Null

$d.Content.Find.Execute(
    " AI's capabilities in data analysis, pattern recognition, and predictive modeling are revolutionizing disease diagnosis, treatment planning, drug discovery, personalized medicine, and epidemic surveillance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Embark on a journey to unravel the enigmatic tapestry of matter, exploring the fundamental principles that govern the behavior of elements and compounds", 2) | Out-Null

$r = $d.Content
$r.Find.Execute(
    " AI-powered systems analyze vast volumes of patient data, identifying complex patterns and correlations that escape human cognition",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Discover the secrets hidden within chemical reactions, witnessing the transformation of substances and the release of energy", 2) | Out-Null
$r.InsertAfter(". Chemistry is not just a collection of abstract concepts; it has tangible applications in everyday life, from the medicines we take to the food we eat")

# --- Body paragraph (second block, after first <w:br/><w:br/>) ---
$d.Content.Find.Execute(
    "Machine learning algorithms, trained on comprehensive datasets, continuously evolve, enhancing their accuracy and effectiveness",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Delve into the fascinating world of chemical elements, the building blocks of the universe", 2) | Out-Null

$d.Content.Find.Execute(
    " As AI systems learn from each new interaction with data, they become more proficient in detecting subtle anomalies, predicting disease progression, and recommending optimal treatments",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Comprehend the periodic table, a roadmap that organizes elements based on their properties", 2) | Out-Null

$r = $d.Content
$r.Find.Execute(
    " AI's impact extends beyond clinical settings, aiding in administrative tasks, enhancing patient engagement, and streamlining healthcare operations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Uncover the secrets of chemical bonding, the forces that hold atoms together, forming molecules and compounds with unique characteristics", 2) | Out-Null
$r.InsertAfter(". Witness the power of chemical reactions, the processes that transform substances into new substances, often accompanied by the release of energy. Explore the intricate web of chemical interactions that occur in living organisms, revealing the intricate mechanisms that govern life")

# --- Body paragraph (third block, after second <w:br/><w:br/>) ---
$d.Content.Find.Execute(
    "AI's integration with medical devices, such as biosensors and wearables, enables continuous monitoring of vital parameters, allowing real-time detection of health issues and immediate intervention",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Chemistry is a dynamic and ever-evolving field, constantly pushing the boundaries of scientific understanding", 2) | Out-Null

$d.Content.Find.Execute(
    " AI-powered telemedicine platforms facilitate remote consultations, expanding access to healthcare for underserved communities and individuals with limited mobility",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " From the synthesis of new materials to the development of innovative drugs, chemistry plays a vital role in shaping the future", 2) | Out-Null

$d.Content.Find.Execute(
    " AI-driven virtual assistants empower patients to actively participate in their healthcare management, providing personalized guidance and support",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Embrace the challenges and rewards of this enigmatic science, expanding your knowledge and gaining a deeper appreciation for the world around you", 2) | Out-Null

# --- Summary paragraph ---
$d.Content.Find.Execute(
    "The fusion of AI and medicine ushers in a new era of healthcare characterized by more precise diagnosis, personalized therapies, and efficient medical operations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In this essay, we delve into the captivating world of chemistry, unveiling the fundamental principles that govern the behavior of matter and the transformations that substances undergo", 2) | Out-Null

$d.Content.Find.Execute(
    " AI's ability to unravel complex data patterns empowers healthcare professionals to make informed decisions, leading to improved patient outcomes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " We explore the periodic table, chemical bonding, chemical reactions, and the intricate web of interactions that occur in living organisms", 2) | Out-Null

$r = $d.Content
$r.Find.Execute(
    " As AI continues to evolve, its potential to transform medicine is limitless, promising a future where healthcare is more proactive, accessible, and effective",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Furthermore, we recognize the practical applications of chemistry in everyday life and appreciate its role in shaping the future through the development of new materials and innovative drugs", 2) | Out-Null
$r.InsertAfter(". Chemistry is not just a subject; it is a lens through which we can comprehend the world and harness its potential to improve lives")

# --- Trailing empty paragraph after the summary ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
